$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ROADM_Type"
$ws.Range("D2:D16").Value = "Directionless"

$ws.Columns.Item(4).ColumnWidth = 13.3

$ws.Range("I5").Select()
